$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1864.2222
$ws.Range("I19").Value = 1388.1111
$ws.Range("J19").Value = 2340.3333
$ws.Range("K19").Value = 1388.1111
$ws.Range("L19").Value = 2340.3333
$ws.Range("M19").Value = -1213.1111
$ws.Range("N19").Value = -2690.3333

$ws.Range("H28").Value = 1148.7333
$ws.Range("I28").Value = 923.3077
$ws.Range("K28").Value = 923.3077
$ws.Range("M28").Value = -438.3077

$ws.Range("H33").Value = 29412942
$ws.Range("I33").Value = 35715524
$ws.Range("K33").Value = 35715524
$ws.Range("M33").Value = -35715295

$ws.Range("H43").Value = 206972.34
$ws.Range("I43").Value = 1622.7273
$ws.Range("J43").Value = 457955.22
$ws.Range("K43").Value = 1622.7273
$ws.Range("L43").Value = 457955.22
$ws.Range("M43").Value = -1553.7273
$ws.Range("N43").Value = -458093.22

$ws.Range("H99").Value = 349
$ws.Range("I99").Value = 349
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1047
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 451
$ws.Range("N99").ClearContents()

$ws.Range("H135").Value = 294716.1
$ws.Range("I135").Value = 323198.3
$ws.Range("J135").Value = 399.66666
$ws.Range("K135").Value = 2908784.7
$ws.Range("L135").Value = 3596.99994
$ws.Range("M135").Value = -2906249.7
$ws.Range("N135").Value = -8666.99994

$ws.Range("H137").Value = 2242.4167
$ws.Range("I137").Value = 2187.762
$ws.Range("J137").Value = 2625
$ws.Range("K137").Value = 6563.286
$ws.Range("L137").Value = 7875
$ws.Range("M137").Value = -4013.286
$ws.Range("N137").Value = -12975

$ws.Range("H138").Value = 4131
$ws.Range("J138").Value = 5319.1694
$ws.Range("L138").Value = 15957.5082
$ws.Range("N138").Value = -26237.5082


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7267.4443
$ws.Range("I45").Value = 4901.1665
$ws.Range("K45").Value = 4901.1665
$ws.Range("M45").Value = -4524.1665

$ws.Range("H57").Value = 4849.1665
$ws.Range("I57").Value = 4849.1665
$ws.Range("K57").Value = 4849.1665
$ws.Range("M57").Value = -4365.1665

$ws.Range("H61").Value = 6310.089
$ws.Range("I61").Value = 3556.0303
$ws.Range("K61").Value = 3556.0303
$ws.Range("M61").Value = -3344.0303

$ws.Range("H74").Value = 29530.027
$ws.Range("I74").Value = 39460.31
$ws.Range("J74").Value = 3711.3
$ws.Range("K74").Value = 39460.31
$ws.Range("L74").Value = 3711.3
$ws.Range("M74").Value = -38586.31
$ws.Range("N74").Value = -5459.3

$ws.Range("H77").Value = 29530.027
$ws.Range("I77").Value = 39460.31
$ws.Range("J77").Value = 3711.3
$ws.Range("K77").Value = 197301.55
$ws.Range("L77").Value = 18556.5
$ws.Range("M77").Value = -192933.55
$ws.Range("N77").Value = -27292.5

$ws.Range("H126").Value = 5265.4443
$ws.Range("I126").Value = 5265.4443
$ws.Range("K126").Value = 15796.3329
$ws.Range("M126").Value = -13326.3329

$ws.Range("H132").Value = 5270.6924
$ws.Range("I132").Value = 1710.7
$ws.Range("K132").Value = 5132.1
$ws.Range("M132").Value = -2602.1

$ws.Range("H136").Value = 6310.089
$ws.Range("I136").Value = 3556.0303
$ws.Range("K136").Value = 10668.0909
$ws.Range("M136").Value = -8118.090899999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 30000
$ws.Range("K82").Value = 30000
$ws.Range("M82").Value = -29617

$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 30000
$ws.Range("K85").Value = 30000
$ws.Range("M85").Value = -28674

$ws.Range("H107").Value = 80360960
$ws.Range("J107").Value = 6166.875
$ws.Range("L107").Value = 6166.875
$ws.Range("N107").Value = -10006.875

$ws.Range("H134").Value = 6318.595
$ws.Range("I134").Value = 2135.7778
$ws.Range("J134").Value = 9455.708000000001
$ws.Range("K134").Value = 6407.3334
$ws.Range("L134").Value = 28367.124
$ws.Range("M134").Value = -3872.3334
$ws.Range("N134").Value = -33437.124


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7250394.5
$ws.Range("I58").Value = 10639664
$ws.Range("K58").Value = 10639664
$ws.Range("M58").Value = -10639461

$ws.Range("H76").Value = 4957.143
$ws.Range("I76").Value = 4957.143
$ws.Range("K76").Value = 4957.143
$ws.Range("M76").Value = -4642.143

$ws.Range("H79").Value = 4957.143
$ws.Range("I79").Value = 4957.143
$ws.Range("K79").Value = 4957.143
$ws.Range("M79").Value = -3865.143

$ws.Range("H132").Value = 3688.8857
$ws.Range("I132").Value = 1583.3529
$ws.Range("K132").Value = 4750.0587
$ws.Range("M132").Value = -2220.0587

$ws.Range("H134").Value = 4241.6616
$ws.Range("I134").Value = 1763.5
$ws.Range("K134").Value = 5290.5
$ws.Range("M134").Value = -2755.5

$ws.Range("H136").Value = 7250394.5
$ws.Range("I136").Value = 10639664
$ws.Range("K136").Value = 31918992
$ws.Range("M136").Value = -31916442


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 12821379
$ws.Range("I14").Value = 12821379
$ws.Range("K14").Value = 38464137
$ws.Range("M14").Value = -38463964

$ws.Range("H92").Value = 6412068
$ws.Range("J92").Value = 6994810.5
$ws.Range("L92").Value = 20984431.5
$ws.Range("N92").Value = -20986927.5

$ws.Range("H98").Value = 2301.0715
$ws.Range("J98").Value = 2417.375
$ws.Range("L98").Value = 7252.125
$ws.Range("N98").Value = -10248.125


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11099.5
$ws.Range("I70").Value = 9798.5
$ws.Range("K70").Value = 9798.5
$ws.Range("M70").Value = -9528.5

$ws.Range("H73").Value = 11099.5
$ws.Range("I73").Value = 9798.5
$ws.Range("K73").Value = 9798.5
$ws.Range("M73").Value = -8862.5

$ws.Range("H122").Value = 3019766.5
$ws.Range("I122").Value = 3813231.5
$ws.Range("K122").Value = 11439694.5
$ws.Range("M122").Value = -11437244.5

$ws.Range("H123").Value = 29974.5
$ws.Range("J123").Value = 29974.5
$ws.Range("L123").Value = 29974.5
$ws.Range("N123").Value = -34874.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 13892150
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 15876385
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 15876385
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -15876761

$ws.Range("H82").Value = 941703.1
$ws.Range("I82").Value = 2819921.8
$ws.Range("K82").Value = 2819921.8
$ws.Range("M82").Value = -2819560.8

$ws.Range("H85").Value = 941703.1
$ws.Range("I85").Value = 2819921.8
$ws.Range("K85").Value = 2819921.8
$ws.Range("M85").Value = -2818673.8

$ws.Range("H122").Value = 4867.2095
$ws.Range("I122").Value = 3616.923
$ws.Range("K122").Value = 10850.769
$ws.Range("M122").Value = -8400.769


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 39589.145
$ws.Range("I62").Value = 104108.8
$ws.Range("J62").Value = 3744.889
$ws.Range("K62").Value = 104108.8
$ws.Range("L62").Value = 3744.889
$ws.Range("M62").Value = -103484.8
$ws.Range("N62").Value = -4992.889

$ws.Range("H65").Value = 39589.145
$ws.Range("I65").Value = 104108.8
$ws.Range("J65").Value = 3744.889
$ws.Range("K65").Value = 520544
$ws.Range("L65").Value = 18724.445
$ws.Range("M65").Value = -517424
$ws.Range("N65").Value = -24964.445

$ws.Range("H136").Value = 23259506
$ws.Range("I136").Value = 35715064
$ws.Range("K136").Value = 107145192
$ws.Range("M136").Value = -107142642

